$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1347.204755070514
$ws.Range("C2").Value = 15729.38127999486
$ws.Range("D2").Value = 335.0348135621319
